{"js": "// Replace each arithmetic-expression answer in the worksheet table with its\n// corrected value. The table is 20 rows x 5 columns; cells are visited in\n// row-major order (left-to-right, top-to-bottom) to match the order the\n// expressions appear in the document, and each old value is verified before\n// being overwritten so a layout mismatch fails loudly instead of silently\n// mis-writing a cell.\nconst replacements = [\n  [\"90-48=42\", \"45+41=86\"], [\"15+63=78\", \"75+21=96\"], [\"41+24=65\", \"95-57=38\"], [\"86+1=87\", \"56-29=27\"], [\"40-4=36\", \"58+28=86\"],\n  [\"32+0=32\", \"11+27=38\"], [\"37-21=16\", \"34+7=41\"], [\"64-48=16\", \"51+11=62\"], [\"15-8=7\", \"90-31=59\"], [\"53-36=17\", \"26+39=65\"],\n  [\"84+14=98\", \"34+47=81\"], [\"44-34=10\", \"0+52=52\"], [\"62-8=54\", \"40+17=57\"], [\"39-39=0\", \"21+38=59\"], [\"68-66=2\", \"73-33=40\"],\n  [\"91-60=31\", \"71-52=19\"], [\"17+53=70\", \"85-74=11\"], [\"69-62=7\", \"14+79=93\"], [\"13+41=54\", \"58+38=96\"], [\"25-1=24\", \"65-56=9\"],\n  [\"19-3=16\", \"72-36=36\"], [\"8+57=65\", \"40+12=52\"], [\"91-36=55\", \"45+52=97\"], [\"76-67=9\", \"8+9=17\"], [\"71-5=66\", \"44-16=28\"],\n  [\"24-13=11\", \"72-60=12\"], [\"32-21=11\", \"42+8=50\"], [\"79-56=23\", \"54-13=41\"], [\"47+46=93\", \"88-77=11\"], [\"0+38=38\", \"53+26=79\"],\n  [\"22+71=93\", \"98-14=84\"], [\"73+0=73\", \"55+41=96\"], [\"27+49=76\", \"38+3=41\"], [\"86-76=10\", \"96-9=87\"], [\"84-79=5\", \"12+74=86\"],\n  [\"56-33=23\", \"54+4=58\"], [\"98-56=42\", \"60-52=8\"], [\"85-24=61\", \"14+45=59\"], [\"94-28=66\", \"2+55=57\"], [\"32+53=85\", \"43+46=89\"],\n  [\"86-39=47\", \"64+25=89\"], [\"39+33=72\", \"59-17=42\"], [\"21+59=80\", \"53-28=25\"], [\"38-24=14\", \"48-33=15\"], [\"53+14=67\", \"52-31=21\"],\n  [\"47-19=28\", \"64+16=80\"], [\"64+13=77\", \"78-21=57\"], [\"57+38=95\", \"39-4=35\"], [\"53-47=6\", \"16+12=28\"], [\"92-82=10\", \"36+57=93\"],\n  [\"89-89=0\", \"77-32=45\"], [\"30+37=67\", \"42+52=94\"], [\"41+8=49\", \"83-52=31\"], [\"84-12=72\", \"90+6=96\"], [\"68+6=74\", \"93-25=68\"],\n  [\"16+81=97\", \"19+31=50\"], [\"22+16=38\", \"59+7=66\"], [\"57+23=80\", \"82-26=56\"], [\"29+0=29\", \"35+49=84\"], [\"78-35=43\", \"6+12=18\"],\n  [\"18+5=23\", \"26+8=34\"], [\"13+40=53\", \"1+90=91\"], [\"88-55=33\", \"10+48=58\"], [\"54-26=28\", \"9+38=47\"], [\"34-7=27\", \"60-46=14\"],\n  [\"77-39=38\", \"56+18=74\"], [\"60-21=39\", \"76-40=36\"], [\"21-20=1\", \"56+2=58\"], [\"78-71=7\", \"25-2=23\"], [\"45-26=19\", \"50-29=21\"],\n  [\"56-25=31\", \"24+56=80\"], [\"34+20=54\", \"8+14=22\"], [\"23+44=67\", \"70+10=80\"], [\"94-83=11\", \"31+60=91\"], [\"54-8=46\", \"55-20=35\"],\n  [\"49-42=7\", \"76-23=53\"], [\"32+54=86\", \"19-1=18\"], [\"91-31=60\", \"92-87=5\"], [\"13+45=58\", \"13+72=85\"], [\"61+9=70\", \"67-32=35\"],\n  [\"90-86=4\", \"80-10=70\"], [\"60+10=70\", \"31-4=27\"], [\"62-17=45\", \"38+61=99\"], [\"44+0=44\", \"56+25=81\"], [\"2+54=56\", \"84-80=4\"],\n  [\"29-14=15\", \"13+79=92\"], [\"40+33=73\", \"90-1=89\"], [\"82-0=82\", \"39+27=66\"], [\"89-33=56\", \"54+2=56\"], [\"79-36=43\", \"97-38=59\"],\n  [\"95-39=56\", \"20+72=92\"], [\"59-22=37\", \"90-48=42\"], [\"39+46=85\", \"84-27=57\"], [\"8+0=8\", \"70-44=26\"], [\"23+61=84\", \"41+38=79\"],\n  [\"89-22=67\", \"97-38=59\"], [\"76-5=71\", \"6+83=89\"], [\"17+5=22\", \"56+33=89\"], [\"15+42=57\", \"26+4=30\"], [\"31+20=51\", \"13+44=57\"],\n];\n\nconst table = context.document.body.tables.getFirstOrNullObject();\ntable.load(\"rowCount\");\nawait context.sync();\n\nif (table.isNullObject) {\n  throw new Error(\"Expected a table in the document body, but none was found.\");\n}\n\nconst columnCount = 5;\ntable.load(\"values\");\nawait context.sync();\n\nconst totalCells = table.rowCount * columnCount;\nif (totalCells !== replacements.length) {\n  throw new Error(\n    `Table shape (${table.rowCount}x${columnCount}=${totalCells}) does not match the expected ${replacements.length} replacements.`\n  );\n}\n\nlet idx = 0;\nfor (let r = 0; r < table.rowCount; r++) {\n  for (let c = 0; c < columnCount; c++) {\n    const [oldText, newText] = replacements[idx];\n    const currentText = table.values[r][c];\n    if (currentText !== oldText) {\n      throw new Error(\n        `Cell (${r},${c}) expected \"${oldText}\" but found \"${currentText}\"; aborting to avoid corrupting unrelated cells.`\n      );\n    }\n    const cell = table.getCell(r, c);\n    cell.value = newText;\n    idx++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each arithmetic-expression answer in the worksheet table with its\n# corrected value. The table is 20 rows x 5 columns; cells are visited in\n# row-major order (left-to-right, top-to-bottom) to match the order the\n# expressions appear in the document, and each old value is verified before\n# being overwritten so a layout mismatch fails loudly instead of silently\n# mis-writing a cell.\n\n$pairs = @(\n    @(\"90-48=42\", \"45+41=86\"),\n    @(\"15+63=78\", \"75+21=96\"),\n    @(\"41+24=65\", \"95-57=38\"),\n    @(\"86+1=87\", \"56-29=27\"),\n    @(\"40-4=36\", \"58+28=86\"),\n    @(\"32+0=32\", \"11+27=38\"),\n    @(\"37-21=16\", \"34+7=41\"),\n    @(\"64-48=16\", \"51+11=62\"),\n    @(\"15-8=7\", \"90-31=59\"),\n    @(\"53-36=17\", \"26+39=65\"),\n    @(\"84+14=98\", \"34+47=81\"),\n    @(\"44-34=10\", \"0+52=52\"),\n    @(\"62-8=54\", \"40+17=57\"),\n    @(\"39-39=0\", \"21+38=59\"),\n    @(\"68-66=2\", \"73-33=40\"),\n    @(\"91-60=31\", \"71-52=19\"),\n    @(\"17+53=70\", \"85-74=11\"),\n    @(\"69-62=7\", \"14+79=93\"),\n    @(\"13+41=54\", \"58+38=96\"),\n    @(\"25-1=24\", \"65-56=9\"),\n    @(\"19-3=16\", \"72-36=36\"),\n    @(\"8+57=65\", \"40+12=52\"),\n    @(\"91-36=55\", \"45+52=97\"),\n    @(\"76-67=9\", \"8+9=17\"),\n    @(\"71-5=66\", \"44-16=28\"),\n    @(\"24-13=11\", \"72-60=12\"),\n    @(\"32-21=11\", \"42+8=50\"),\n    @(\"79-56=23\", \"54-13=41\"),\n    @(\"47+46=93\", \"88-77=11\"),\n    @(\"0+38=38\", \"53+26=79\"),\n    @(\"22+71=93\", \"98-14=84\"),\n    @(\"73+0=73\", \"55+41=96\"),\n    @(\"27+49=76\", \"38+3=41\"),\n    @(\"86-76=10\", \"96-9=87\"),\n    @(\"84-79=5\", \"12+74=86\"),\n    @(\"56-33=23\", \"54+4=58\"),\n    @(\"98-56=42\", \"60-52=8\"),\n    @(\"85-24=61\", \"14+45=59\"),\n    @(\"94-28=66\", \"2+55=57\"),\n    @(\"32+53=85\", \"43+46=89\"),\n    @(\"86-39=47\", \"64+25=89\"),\n    @(\"39+33=72\", \"59-17=42\"),\n    @(\"21+59=80\", \"53-28=25\"),\n    @(\"38-24=14\", \"48-33=15\"),\n    @(\"53+14=67\", \"52-31=21\"),\n    @(\"47-19=28\", \"64+16=80\"),\n    @(\"64+13=77\", \"78-21=57\"),\n    @(\"57+38=95\", \"39-4=35\"),\n    @(\"53-47=6\", \"16+12=28\"),\n    @(\"92-82=10\", \"36+57=93\"),\n    @(\"89-89=0\", \"77-32=45\"),\n    @(\"30+37=67\", \"42+52=94\"),\n    @(\"41+8=49\", \"83-52=31\"),\n    @(\"84-12=72\", \"90+6=96\"),\n    @(\"68+6=74\", \"93-25=68\"),\n    @(\"16+81=97\", \"19+31=50\"),\n    @(\"22+16=38\", \"59+7=66\"),\n    @(\"57+23=80\", \"82-26=56\"),\n    @(\"29+0=29\", \"35+49=84\"),\n    @(\"78-35=43\", \"6+12=18\"),\n    @(\"18+5=23\", \"26+8=34\"),\n    @(\"13+40=53\", \"1+90=91\"),\n    @(\"88-55=33\", \"10+48=58\"),\n    @(\"54-26=28\", \"9+38=47\"),\n    @(\"34-7=27\", \"60-46=14\"),\n    @(\"77-39=38\", \"56+18=74\"),\n    @(\"60-21=39\", \"76-40=36\"),\n    @(\"21-20=1\", \"56+2=58\"),\n    @(\"78-71=7\", \"25-2=23\"),\n    @(\"45-26=19\", \"50-29=21\"),\n    @(\"56-25=31\", \"24+56=80\"),\n    @(\"34+20=54\", \"8+14=22\"),\n    @(\"23+44=67\", \"70+10=80\"),\n    @(\"94-83=11\", \"31+60=91\"),\n    @(\"54-8=46\", \"55-20=35\"),\n    @(\"49-42=7\", \"76-23=53\"),\n    @(\"32+54=86\", \"19-1=18\"),\n    @(\"91-31=60\", \"92-87=5\"),\n    @(\"13+45=58\", \"13+72=85\"),\n    @(\"61+9=70\", \"67-32=35\"),\n    @(\"90-86=4\", \"80-10=70\"),\n    @(\"60+10=70\", \"31-4=27\"),\n    @(\"62-17=45\", \"38+61=99\"),\n    @(\"44+0=44\", \"56+25=81\"),\n    @(\"2+54=56\", \"84-80=4\"),\n    @(\"29-14=15\", \"13+79=92\"),\n    @(\"40+33=73\", \"90-1=89\"),\n    @(\"82-0=82\", \"39+27=66\"),\n    @(\"89-33=56\", \"54+2=56\"),\n    @(\"79-36=43\", \"97-38=59\"),\n    @(\"95-39=56\", \"20+72=92\"),\n    @(\"59-22=37\", \"90-48=42\"),\n    @(\"39+46=85\", \"84-27=57\"),\n    @(\"8+0=8\", \"70-44=26\"),\n    @(\"23+61=84\", \"41+38=79\"),\n    @(\"89-22=67\", \"97-38=59\"),\n    @(\"76-5=71\", \"6+83=89\"),\n    @(\"17+5=22\", \"56+33=89\"),\n    @(\"15+42=57\", \"26+4=30\"),\n    @(\"31+20=51\", \"13+44=57\")\n)\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n$rowCount = $tbl.Rows.Count\n$colCount = $tbl.Columns.Count\n\nif (($rowCount * $colCount) -ne $pairs.Count) {\n    throw \"Table shape ($rowCount x $colCount = $($rowCount * $colCount)) does not match the expected $($pairs.Count) replacements.\"\n}\n\n$idx = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n    for ($c = 1; $c -le $colCount; $c++) {\n        $cell = $tbl.Cell($r, $c)\n        $oldText = $pairs[$idx][0]\n        $newText = $pairs[$idx][1]\n\n        # Cell text includes the trailing cell-mark (CR + BEL); strip it before comparing.\n        $currentText = $cell.Range.Text.TrimEnd([char]13, [char]7)\n\n        if ($currentText -ne $oldText) {\n            throw \"Cell ($r,$c) expected '$oldText' but found '$currentText'; aborting to avoid corrupting unrelated cells.\"\n        }\n\n        $cell.Range.Text = $newText\n        $idx++\n    }\n}\n"}
